# Apply the GDD.docx content edits described by the commit.
$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $find"
    }
}

$rsq = [char]0x2019  # right single quotation mark

# --- Paragraph 1 ("About: ...") -----------------------------------------
$find1 = "Balloon fight/joust style game where colour of the screen changes changes and the physics change with it. The player is in a room with an open pit at the bottom. Falling in the pit is instant death. Enemies will constantly spawn and attack the player. They will run at the player without fear and try to push the player off into the pit. When the player hit${rsq}s the pit it is game over and the player get${rsq}s the option to restart or play again. The goal is to score as many points as possible without dying and then beat the high score. "
$replace1 = "Colour Splatter eXtreme is single screen score attack game inspired by arcade classics Balloon fight and Joust. The defining gimmick is that every 20 seconds the screen changes colour. When the colour of the screen changes the physics of the game world change with it. "
Replace-Text $find1 $replace1

# --- Paragraph 2 ("Every 20 seconds ...") now becomes the player paragraph
$find2 = "Every 20 seconds the colour of the light of screen will change. With that changes some amount of physics will change according to the list below."
$replace2 = "The player is in a room with an open pit at the bottom. Falling in the pit is instant death. Enemies will constantly spawn and attack the player. They will run at the player without fear and try to push the player into the pit. When the player falls into the pit it is game over and the player has the option to restart or play again. The goal is to score as many points as possible by killing enemies, without dying and beat the high score. "
Replace-Text $find2 $replace2

# --- Colours list -----------------------------------------------------
Replace-Text "1 - Blue: mass of all object becomes lower" "1 - Blue: Stronger Gravity (added mass). Objects will feel heavier and will slow down more quickly from momentum"
Replace-Text "2 - Red: velocity grows higher" "2 - Red: Velocity grows higher. Object will move faster within the world."
Replace-Text "3 - Green: Everything becomes 3x bigger" "3 - Green: All Objects become 3x bigger."
Replace-Text "4 - Orange: Bouncy" "4 - Orange: The ground becomes bouncy and hard to stand still on."
Replace-Text "5 - Purple: Objects stick together (delay on jump)" "5 - Purple: The ground becomes sticky (delay on jump) off ground"

# Collapse "8 - " / "Yellow: Random (between 1 and 6" / ")" into one run
# (text is unchanged, but the source runs get merged into a single run).
Replace-Text "8 - Yellow: Random (between 1 and 6)" "8 - Yellow: Random (between 1 and 6)"

# --- Move the automatic "_GoBack" bookmark to follow the last edit -----
# (mirrors Word's own behaviour of re-stamping _GoBack at the most
#  recently-typed location -- it previously sat around "Angry Tree").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$bmRng = $d.Content
$bmRng.Find.Execute("heavier", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmRng.Collapse(0) # wdCollapseEnd
$d.Bookmarks.Add("_GoBack", $bmRng) | Out-Null

# --- Player character description: "Splat" -> "Splater" ----------------
Replace-Text "Colour Splat eXtreme competition." "Colour Splater eXtreme competition."

# --- Notes: add new final sentence --------------------------------------
Replace-Text "Enemies will spawn every 3-7 seconds assuming the cap of enemies has not been hit." "Enemies will spawn every 3-7 seconds assuming the cap of enemies has not been hit.`rColour changes every 20 seconds"

Write-Output "done"
